$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5884189605712891
$ws.Range("B1").Value = 2.064315557479858
$ws.Range("D1").Value = 2.673370122909546
$ws.Range("E1").Value = 1.08525812625885
